# Update (Analyze PO & Forecast)
$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": update MyForecast (column D) numeric values ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("D2").Value = 44
$wsForecast.Range("D3").Value = 37
$wsForecast.Range("D4").Value = 36
$wsForecast.Range("D5").Value = 42
$wsForecast.Range("D6").Value = 47
$wsForecast.Range("D7").Value = 47
$wsForecast.Range("D8").Value = 46
$wsForecast.Range("D9").Value = 43
$wsForecast.Range("D10").Value = 41
$wsForecast.Range("D11").Value = 43
$wsForecast.Range("D12").Value = 42
$wsForecast.Range("D13").Value = 42
$wsForecast.Range("D14").Value = 41

# --- Sheet "Summary": update derived metrics (column B) ---
# These cells are stored as text in the workbook, so force text format
# before assigning the new values to preserve their original text type.
$wsSummary = $wb.Worksheets.Item("Summary")

$rngSummaryB = $wsSummary.Range("B9:B15")
$rngSummaryB.NumberFormat = "@"

$wsSummary.Range("B9").Value = "669"
$wsSummary.Range("B10").Value = "341"
$wsSummary.Range("B11").Value = "159"
$wsSummary.Range("B12").Value = "47"
$wsSummary.Range("B13").Value = "2025-02-23"
$wsSummary.Range("B14").Value = "36"
$wsSummary.Range("B15").Value = "2025-02-09"
